# Update cryptos list with latest price/volume data (refresh run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: Price values in column D are plain text (not real numbers, since some
# use dotted thousands-separators like "69.382.49"). We prefix numeric-looking
# values with a leading single quote so Excel stores them as text, matching
# the workbook's inlineStr string cells instead of converting them to Double.

$ws.Range("D2").Value = '''69.382.49'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '''3.443.90'
$ws.Range("E3").Value = '  +2.28%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''584.26'
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").Value = '''179.23'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '''0.199'
$ws.Range("E9").Value = '  +7.26%  '
$ws.Range("D10").Value = '''0.587'
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("D11").Value = '''48.53'
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").Value = '''684.77'
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").Value = '''3.986.65'
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = '''69.451.55'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '''3.437.71'
$ws.Range("E17").Value = '  +1.12%  '
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = '''17.88'
$ws.Range("E19").Value = '  +1.54%  '
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").Value = '''0.910'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").Value = '''5.37'
$ws.Range("E22").Value = '  -2.37%  '
$ws.Range("D23").Value = '''17.04'
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("D24").Value = '''101.14'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("D26").Value = '''2.70'
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("D28").Value = '''33.68'
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("D29").Value = '''8.80'
$ws.Range("E29").Value = '  +2.14%  '
$ws.Range("D30").Value = '''6.90'
$ws.Range("E30").Value = '  -1.81%  '
$ws.Range("E31").Value = '  +7.36%  '
$ws.Range("D32").Value = '''563.44'
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("E33").Value = '  -0.65%  '
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("D35").Value = '''58.19'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").Value = '''3.626.82'
$ws.Range("E37").Value = '  -2.35%  '
$ws.Range("E38").Value = '  -1.14%  '
$ws.Range("D39").Value = '''35.13'
$ws.Range("E39").Value = '  +0.95%  '
$ws.Range("D40").Value = '''0.0₃0742'
$ws.Range("E40").Value = '  +9.08%  '
$ws.Range("E41").Value = '  +2.67%  '
$ws.Range("D42").Value = '''2.70'
$ws.Range("E42").Value = '  +2.28%  '
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").Value = '''3.35'
$ws.Range("E43").Value = '  +2.86%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '''0.0425'
$ws.Range("E44").Value = '  +2.02%  '
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("E48").Value = '  +4.86%  '
$ws.Range("D49").Value = '''0.999'
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("D50").Value = '''131.39'
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").Value = '''2.69'
$ws.Range("E51").Value = '  +1.94%  '
